$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("/NewDataSet/Table[2]/Location", "[A-Z a-z].*"),
    @("/NewDataSet/Table[2]/PostCode", "[A-Z a-z 0-9].*"),
    @("/NewDataSet/Table[3]/Location", "[A-Z a-z].*"),
    @("/NewDataSet/Table[3]/PostCode", "[A-Z a-z 0-9].*"),
    @("/NewDataSet/Table[4]/Location", "[A-Z a-z].*"),
    @("/NewDataSet/Table[4]/PostCode", "[A-Z a-z 0-9].*")
)

$r = 4
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $r++
}

$ws.Range("D9").Select()
